$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
# A3: 111701910 -> 111701829
$ws.Range("A3").Value = 111701829
# P3/Q3/R3: location changes to "Myrövägen öster" coordinates
$ws.Range("P3").Value = "Myrövägen öster (Myrövägen öster), Nrk"
$ws.Range("Q3").Value = 516894.5773385105
$ws.Range("R3").Value = 6574639.474785783

# --- Row 5 ---
# A5: 111701829 -> 111702271
$ws.Range("A5").Value = 111702271
# B5: 90687 -> 90709
$ws.Range("B5").Value = 90709
# D5: LC -> NT
$ws.Range("D5").Value = "NT"
# E5: 5964 -> 5448
$ws.Range("E5").Value = 5448
# F5: Fjällig taggsvamp s.str. -> Svartvit taggsvamp
$ws.Range("F5").Value = "Svartvit taggsvamp"
# G5: Sarcodon imbricatus s.str. -> Phellodon connatus
$ws.Range("G5").Value = "Phellodon connatus"
# H5: (L.:Fr.) P.Karst. -> (Schultz) nom.prov
$ws.Range("H5").Value = "(Schultz) nom.prov"
# I5: (empty) -> "3" (must stay text, not become the number 3)
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "3"
# J5: (empty) -> "fruktkroppar"
$ws.Range("J5").Value = "fruktkroppar"
# P5/Q5/R5: location changes to "Kyrkogården" coordinates
$ws.Range("P5").Value = "Kyrkogården (Kyrkogården), Nrk"
$ws.Range("Q5").Value = 516923.6092008445
$ws.Range("R5").Value = 6574666.663922376

# --- Row 6 ---
# A6: 111702271 -> 111701910
$ws.Range("A6").Value = 111701910
# B6: 90709 -> 90687
$ws.Range("B6").Value = 90687
# D6: NT -> LC
$ws.Range("D6").Value = "LC"
# E6: 5448 -> 5964
$ws.Range("E6").Value = 5964
# F6: Svartvit taggsvamp -> Fjällig taggsvamp s.str.
$ws.Range("F6").Value = "Fjällig taggsvamp s.str."
# G6: Phellodon connatus -> Sarcodon imbricatus s.str.
$ws.Range("G6").Value = "Sarcodon imbricatus s.str."
# H6: (Schultz) nom.prov -> (L.:Fr.) P.Karst.
$ws.Range("H6").Value = "(L.:Fr.) P.Karst."
# I6: "3" -> (empty)
$ws.Range("I6").Value = ""
# J6: "fruktkroppar" -> (empty)
$ws.Range("J6").Value = ""
# P6 unchanged ("Kyrkogården (Kyrkogården), Nrk")
# Q6/R6: coordinates revert to original row3 coordinates
$ws.Range("Q6").Value = 516978.9846792166
$ws.Range("R6").Value = 6574635.767148005
